# TimeSheets.xlsx update
# - Refresh of the "project_fynbus_time_tracking" CSV data connection pulled in
#   new rows (re-labelled "(1)" since it's a re-downloaded copy of the CSV) on
#   both Sheet1 and Sheet2, inserting 2 rows of new data before the old
#   "Tests for Ordering Offers" block and appending 5 more rows after it.
# - The corresponding sheet-scoped defined name is renamed/resized to match.

$wb = $excel.ActiveWorkbook

function Update-TimesheetSheet($ws) {
    # Insert 2 new rows right before the existing row 68 ("Tests for Ordering
    # Offers" block), pushing the old rows 68-71 down to 70-73 unchanged.
    $ws.Rows("68:69").Insert()

    # New row 68
    $ws.Range("A68").Value = 42816
    $ws.Range("B68").Value = "-"
    $ws.Range("C68").Value = "Project Fynbus"
    $ws.Range("D68").Value = "Fulfill Offers Test"
    $ws.Range("E68").Value = 1.4861
    $ws.Range("F68").Value = "Roxana"
    $ws.Range("G68").Value = "Ion"

    # New row 69
    $ws.Range("A69").Value = 42816
    $ws.Range("B69").Value = "-"
    $ws.Range("C69").Value = "Project Fynbus"
    $ws.Range("D69").Value = "Fulfill Offers Test"
    $ws.Range("E69").Value = 1.5261
    $ws.Range("F69").Value = "Hedviga"
    $ws.Range("G69").Value = "Arta Gerina"

    # Rows 70-73 already hold the old rows 68-71 data (shifted down by the
    # insert) and require no further changes.

    # Apply the date number format (matching column A elsewhere) to the newly
    # appended rows 74-78 before filling them in.
    $ws.Range("A67").Copy()
    $ws.Range("A74:A78").PasteSpecial(-4122)

    # New row 74
    $ws.Range("A74").Value = 42817
    $ws.Range("B74").Value = "-"
    $ws.Range("C74").Value = "Project Fynbus"
    $ws.Range("D74").Value = "Fulfill Offers Test"
    $ws.Range("E74").Value = 0.3167
    $ws.Range("F74").Value = "Roxana"
    $ws.Range("G74").Value = "Ion"

    # New row 75
    $ws.Range("A75").Value = 42817
    $ws.Range("B75").Value = "-"
    $ws.Range("C75").Value = "Project Fynbus"
    $ws.Range("D75").Value = "Update ClassDiagram"
    $ws.Range("E75").Value = 0.0517
    $ws.Range("F75").Value = "Jonas"
    $ws.Range("G75").Value = "Laursen"

    # New row 76
    $ws.Range("A76").Value = 42817
    $ws.Range("B76").Value = "-"
    $ws.Range("C76").Value = "Project Fynbus"
    $ws.Range("D76").Value = "Test for Routes"
    $ws.Range("E76").Value = 0.3761
    $ws.Range("F76").Value = "Jonas"
    $ws.Range("G76").Value = "Laursen"

    # New row 77
    $ws.Range("A77").Value = 42817
    $ws.Range("B77").Value = "-"
    $ws.Range("C77").Value = "Project Fynbus"
    $ws.Range("D77").Value = "Tests for Ordering Offers"
    $ws.Range("E77").Value = 0.4067
    $ws.Range("F77").Value = "Matthew"
    $ws.Range("G77").Value = "Peterson"

    # New row 78
    $ws.Range("A78").Value = 42817
    $ws.Range("B78").Value = "-"
    $ws.Range("C78").Value = "Project Fynbus"
    $ws.Range("D78").Value = "Tests for Ordering Offers"
    $ws.Range("E78").Value = 0.4036
    $ws.Range("F78").Value = "Jonas"
    $ws.Range("G78").Value = "Laursen"
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

Update-TimesheetSheet $ws1
Update-TimesheetSheet $ws2

# Re-point the Sheet2-scoped defined name at the grown range and rename it
# (matches the "_1" suffix the CSV re-import gave it).
$wb.Names.Item("Sheet2!project_fynbus_time_tracking").Delete()
$ws2.Names.Add("project_fynbus_time_tracking_1", "=Sheet2!`$A`$1:`$G`$78")

# Update the view state: scrolled down near the bottom of the sheet with the
# newly-added rows selected.
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws1.Range("A58:G78").Select()

$ws2.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 52
$win2.ScrollColumn = 1
$ws2.Range("A58:G78").Select()

$ws1.Activate()
